$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.06503766666666667
$ws.Range("H2").Value = 0.195113
$ws.Range("I2").Value = 0.009380539125347769
$ws.Range("J2").Value = 0.009380539125347769
$ws.Range("M2").Value = 281.0920463333333
$ws.Range("N2").Value = 843.2761389999999
$ws.Range("O2").Value = 0.8291026083535286
$ws.Range("P2").Value = 0.8291026083535286
$ws.Range("Q2").Value = 18.28157081207856
$ws.Range("R2").Value = 164.534137308707
$ws.Range("S2").Value = 0.007777429456588163
$ws.Range("T2").Value = 0.007777429456588163
$ws.Range("G3").Value = 0.06503766666666667
$ws.Range("H3").Value = 0.195113
$ws.Range("I3").Value = 0.009380539125347769
$ws.Range("J3").Value = 0.009380539125347769
$ws.Range("O3").Value = 0.001324719879221983
$ws.Range("P3").Value = 0.001324719879221983
$ws.Range("Q3").Value = 0.02920984692866667
$ws.Range("R3").Value = 0.262888622358
$ws.Range("S3").Value = 0.00001242658665716779
$ws.Range("T3").Value = 0.00001242658665716779
$ws.Range("G4").Value = 0.06503766666666667
$ws.Range("H4").Value = 0.195113
$ws.Range("I4").Value = 0.009380539125347769
$ws.Range("J4").Value = 0.009380539125347769
$ws.Range("M4").Value = 4.452417
$ws.Range("N4").Value = 13.357251
$ws.Range("O4").Value = 0.01313274635953239
$ws.Range("P4").Value = 0.01313274635953239
$ws.Range("Q4").Value = 0.289574812707
$ws.Range("R4").Value = 2.606173314363
$ws.Range("S4").Value = 0.000123192241048862
$ws.Range("T4").Value = 0.000123192241048862
$ws.Range("G5").Value = 0.06503766666666667
$ws.Range("H5").Value = 0.195113
$ws.Range("I5").Value = 0.009380539125347769
$ws.Range("J5").Value = 0.009380539125347769
$ws.Range("M5").Value = 53.03808999999999
$ws.Range("N5").Value = 159.11427
$ws.Range("O5").Value = 0.156439925407717
$ws.Range("P5").Value = 0.156439925407717
$ws.Range("Q5").Value = 3.449473618056667
$ws.Range("R5").Value = 31.04526256251
$ws.Range("S5").Value = 0.001467490841053576
$ws.Range("T5").Value = 0.001467490841053576
$ws.Range("G6").Value = 4.613664666666667
$ws.Range("I6").Value = 0.6654399540302477
$ws.Range("J6").Value = 0.6654399540302477
$ws.Range("M6").Value = 281.0920463333333
$ws.Range("N6").Value = 843.2761389999999
$ws.Range("O6").Value = 0.8291026083535286
$ws.Range("P6").Value = 0.8291026083535286
$ws.Range("Q6").Value = 1296.86444224913
$ws.Range("R6").Value = 11671.77998024217
$ws.Range("S6").Value = 0.5517180015891305
$ws.Range("T6").Value = 0.5517180015891305
$ws.Range("G7").Value = 4.613664666666667
$ws.Range("I7").Value = 0.6654399540302477
$ws.Range("J7").Value = 0.6654399540302477
$ws.Range("O7").Value = 0.001324719879221983
$ws.Range("P7").Value = 0.001324719879221983
$ws.Range("S7").Value = 0.0008815215355324319
$ws.Range("T7").Value = 0.0008815215355324319
$ws.Range("G8").Value = 4.613664666666667
$ws.Range("I8").Value = 0.6654399540302477
$ws.Range("J8").Value = 0.6654399540302477
$ws.Range("M8").Value = 4.452417
$ws.Range("N8").Value = 13.357251
$ws.Range("O8").Value = 0.01313274635953239
$ws.Range("P8").Value = 0.01313274635953239
$ws.Range("Q8").Value = 20.541958994166
$ws.Range("R8").Value = 184.877630947494
$ws.Range("S8").Value = 0.008739054133778134
$ws.Range("T8").Value = 0.008739054133778136
$ws.Range("G9").Value = 4.613664666666667
$ws.Range("I9").Value = 0.6654399540302477
$ws.Range("J9").Value = 0.6654399540302477
$ws.Range("M9").Value = 53.03808999999999
$ws.Range("N9").Value = 159.11427
$ws.Range("O9").Value = 0.156439925407717
$ws.Range("P9").Value = 0.156439925407717
$ws.Range("Q9").Value = 244.6999618204866
$ws.Range("R9").Value = 2202.29965638438
$ws.Range("S9").Value = 0.1041013767718066
$ws.Range("T9").Value = 0.1041013767718066
$ws.Range("G10").Value = 2.254552333333333
$ws.Range("H10").Value = 6.763657
$ws.Range("I10").Value = 0.3251795068444046
$ws.Range("J10").Value = 0.3251795068444046
$ws.Range("M10").Value = 281.0920463333333
$ws.Range("N10").Value = 843.2761389999999
$ws.Range("O10").Value = 0.8291026083535286
$ws.Range("P10").Value = 0.8291026083535286
$ws.Range("Q10").Value = 633.7367289422581
$ws.Range("R10").Value = 5703.630560480323
$ws.Range("S10").Value = 0.2696071773078099
$ws.Range("T10").Value = 0.26960717730781
$ws.Range("G11").Value = 2.254552333333333
$ws.Range("H11").Value = 6.763657
$ws.Range("I11").Value = 0.3251795068444046
$ws.Range("J11").Value = 0.3251795068444046
$ws.Range("O11").Value = 0.001324719879221983
$ws.Range("P11").Value = 0.001324719879221983
$ws.Range("Q11").Value = 1.012569053051333
$ws.Range("R11").Value = 9.113121477462
$ws.Range("S11").Value = 0.0004307717570323838
$ws.Range("T11").Value = 0.0004307717570323838
$ws.Range("G12").Value = 2.254552333333333
$ws.Range("H12").Value = 6.763657
$ws.Range("I12").Value = 0.3251795068444046
$ws.Range("J12").Value = 0.3251795068444046
$ws.Range("M12").Value = 4.452417
$ws.Range("N12").Value = 13.357251
$ws.Range("O12").Value = 0.01313274635953239
$ws.Range("P12").Value = 0.01313274635953239
$ws.Range("Q12").Value = 10.038207136323
$ws.Range("R12").Value = 90.343864226907
$ws.Range("S12").Value = 0.004270499984705391
$ws.Range("T12").Value = 0.004270499984705392
$ws.Range("G13").Value = 2.254552333333333
$ws.Range("H13").Value = 6.763657
$ws.Range("I13").Value = 0.3251795068444046
$ws.Range("J13").Value = 0.3251795068444046
$ws.Range("M13").Value = 53.03808999999999
$ws.Range("N13").Value = 159.11427
$ws.Range("O13").Value = 0.156439925407717
$ws.Range("P13").Value = 0.156439925407717
$ws.Range("Q13").Value = 119.5771495650433
$ws.Range("R13").Value = 1076.19434608539
$ws.Range("S13").Value = 0.05087105779485684
$ws.Range("T13").Value = 0.05087105779485686
